$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text changes (shared string content updates) ---
$ws.Range("J4").Value = "All expected passes"
$ws.Range("J13").Value = "All expected passes"

$ws.Range("J36").Value = "passed"
$ws.Range("J37").Value = "unexpected fail"
$ws.Range("J41").Value = "passed"

$ws.Range("J53").Value = "FunnyQT" + [char]0x2019 + "s BX-transformation framework is essentially state-based and not incremental (listening to model changes) and doesn" + [char]0x2019 + "t preserve correspondences between transformation executions, so operations like renaming or moving are resolved by deletion and re-creation which might loose information manually set in the target model (like the birthdays)."

# --- Numeric metric updates ---
$ws.Range("J45").Value = 52
$ws.Range("J46").Value = 264
$ws.Range("J47").Value = 2172

# --- Row height update (row 53 grew taller to fit new text) ---
$ws.Rows("53").RowHeight = 161.2

# --- View: scroll / freeze panes ---
$excel.ActiveWindow.ScrollColumn = 9
$excel.ActiveWindow.ScrollRow = 46
$ws.Range("E1").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("J54").Select()
